# Apply the "push from new laptop" edit to the vacancy workbook.
$wb = $excel.ActiveWorkbook

$wsVacancy = $wb.Worksheets.Item("vacancy")
$wsEdit = $wb.Worksheets.Item("edit")

# --- Sheet "vacancy" (sheet1.xml) ---

# Update site / month text values on rows 2-4.
$wsVacancy.Range("B2").Value = "100011 - apollo"
$wsVacancy.Range("D2").Value = "December"

$wsVacancy.Range("B3").Value = "100010 - apollo"
$wsVacancy.Range("D3").Value = "January"

$wsVacancy.Range("A4").Value = " apollo medical services"
$wsVacancy.Range("B4").Value = "10012 - apollo medical services"
$wsVacancy.Range("D4").Value = "January"

# Drop the now-unused "NoOfClosedvacancies" column data (G), keep column layout.
$wsVacancy.Range("G1:G4").ClearContents()

# Widen column B slightly (closest the engine's width quantization allows to 29.7109375).
$wsVacancy.Columns.Item(2).ColumnWidth = 28.83

# Style row 4: A4 right aligned, B4 wrap text, taller row.
$wsVacancy.Range("A4").HorizontalAlignment = -4152  # xlRight
$wsVacancy.Range("B4").WrapText = $true
$wsVacancy.Rows.Item(4).RowHeight = 30

# --- Sheet "edit" (sheet2.xml) ---

# Drop the now-unused "No Of Closed Vacancies" column data (C), keep column layout.
$wsEdit.Range("C1:C3").ClearContents()

$wsEdit.Range("C5").Select()

$wb.Save()
